$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10
$ws.Range("A10").Value = 111936866
$ws.Range("B10").Value = 89423
$ws.Range("C10").Value = "Ovaliderad"
$ws.Range("D10").Value = "NT"
$ws.Range("E10").Value = 5432
$ws.Range("F10").Value = "Granticka"
$ws.Range("G10").Value = "Porodaedalea chrysoloma"
$ws.Range("H10").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("I10").Value = "'"
$ws.Range("P10").Value = "rörvattsbodarna, Jmt"
$ws.Range("Q10").Value = 448765.5992023234
$ws.Range("R10").Value = 7087416.731054713
$ws.Range("S10").Value = 10
$ws.Range("T10").Value = "Jämtland"
$ws.Range("U10").Value = "Krokom"
$ws.Range("V10").Value = "Jämtland"
$ws.Range("W10").Value = "Hotagen"
$ws.Range("Y10").Value = "'2023-09-05"
$ws.Range("Z10").Value = "00:00"
$ws.Range("AA10").Value = "'2023-09-05"
$ws.Range("AB10").Value = "00:00"
$ws.Range("AD10").Value = $false
$ws.Range("AE10").Value = $false
$ws.Range("AG10").Value = $false
$ws.Range("AT10").Value = "'"
$ws.Range("AW10").Value = "Benny Öwre"
$ws.Range("AX10").Value = "Benny Öwre"
$ws.Range("AY10").Value = "'"

# Row 11
$ws.Range("A11").Value = 111936792
$ws.Range("B11").Value = 90087
$ws.Range("C11").Value = "Ovaliderad"
$ws.Range("D11").Value = "LC"
$ws.Range("E11").Value = 3298
$ws.Range("F11").Value = "Trådticka"
$ws.Range("G11").Value = "Climacocystis borealis"
$ws.Range("H11").Value = "(Fr.) Kotl. & Pouzar"
$ws.Range("I11").Value = "'"
$ws.Range("P11").Value = "rörvattsbodarna, Jmt"
$ws.Range("Q11").Value = 448761.1990147882
$ws.Range("R11").Value = 7087578.827763715
$ws.Range("S11").Value = 10
$ws.Range("T11").Value = "Jämtland"
$ws.Range("U11").Value = "Krokom"
$ws.Range("V11").Value = "Jämtland"
$ws.Range("W11").Value = "Hotagen"
$ws.Range("Y11").Value = "'2023-09-05"
$ws.Range("Z11").Value = "00:00"
$ws.Range("AA11").Value = "'2023-09-05"
$ws.Range("AB11").Value = "00:00"
$ws.Range("AD11").Value = $false
$ws.Range("AE11").Value = $false
$ws.Range("AG11").Value = $false
$ws.Range("AT11").Value = "'"
$ws.Range("AW11").Value = "Benny Öwre"
$ws.Range("AX11").Value = "Benny Öwre"
$ws.Range("AY11").Value = "'"

# Row 12
$ws.Range("A12").Value = 111936867
$ws.Range("B12").Value = 89423
$ws.Range("C12").Value = "Ovaliderad"
$ws.Range("D12").Value = "NT"
$ws.Range("E12").Value = 5432
$ws.Range("F12").Value = "Granticka"
$ws.Range("G12").Value = "Porodaedalea chrysoloma"
$ws.Range("H12").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("I12").Value = "'"
$ws.Range("P12").Value = "rörvattsbodarna, Jmt"
$ws.Range("Q12").Value = 448791.554596175
$ws.Range("R12").Value = 7087386.366048628
$ws.Range("S12").Value = 10
$ws.Range("T12").Value = "Jämtland"
$ws.Range("U12").Value = "Krokom"
$ws.Range("V12").Value = "Jämtland"
$ws.Range("W12").Value = "Hotagen"
$ws.Range("Y12").Value = "'2023-09-05"
$ws.Range("Z12").Value = "00:00"
$ws.Range("AA12").Value = "'2023-09-05"
$ws.Range("AB12").Value = "00:00"
$ws.Range("AD12").Value = $false
$ws.Range("AE12").Value = $false
$ws.Range("AG12").Value = $false
$ws.Range("AT12").Value = "'"
$ws.Range("AW12").Value = "Benny Öwre"
$ws.Range("AX12").Value = "Benny Öwre"
$ws.Range("AY12").Value = "'"

# Row 13
$ws.Range("A13").Value = 111936798
$ws.Range("B13").Value = 56398
$ws.Range("C13").Value = "Ovaliderad"
$ws.Range("D13").Value = "NT"
$ws.Range("E13").Value = 100109
$ws.Range("F13").Value = "Tretåig hackspett"
$ws.Range("G13").Value = "Picoides tridactylus"
$ws.Range("H13").Value = "(Linnaeus, 1758)"
$ws.Range("I13").Value = "'"
$ws.Range("K13").Value = "'"
$ws.Range("L13").Value = "'"
$ws.Range("M13").Value = "'"
$ws.Range("N13").Value = "'"
$ws.Range("P13").Value = "rörvattsbodarna, Jmt"
$ws.Range("Q13").Value = 448923.1251473558
$ws.Range("R13").Value = 7087371.00725084
$ws.Range("S13").Value = 10
$ws.Range("T13").Value = "Jämtland"
$ws.Range("U13").Value = "Krokom"
$ws.Range("V13").Value = "Jämtland"
$ws.Range("W13").Value = "Hotagen"
$ws.Range("Y13").Value = "'2023-09-05"
$ws.Range("Z13").Value = "00:00"
$ws.Range("AA13").Value = "'2023-09-05"
$ws.Range("AB13").Value = "00:00"
$ws.Range("AC13").Value = "ringhack äldre"
$ws.Range("AD13").Value = $false
$ws.Range("AE13").Value = $false
$ws.Range("AG13").Value = $false
$ws.Range("AT13").Value = "'"
$ws.Range("AW13").Value = "Benny Öwre"
$ws.Range("AX13").Value = "Benny Öwre"
$ws.Range("AY13").Value = "'"

# Row 14
$ws.Range("A14").Value = 111936868
$ws.Range("B14").Value = 89423
$ws.Range("C14").Value = "Ovaliderad"
$ws.Range("D14").Value = "NT"
$ws.Range("E14").Value = 5432
$ws.Range("F14").Value = "Granticka"
$ws.Range("G14").Value = "Porodaedalea chrysoloma"
$ws.Range("H14").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("I14").Value = "'"
$ws.Range("P14").Value = "rörvattsbodarna, Jmt"
$ws.Range("Q14").Value = 448988.017639213
$ws.Range("R14").Value = 7087186.778340456
$ws.Range("S14").Value = 10
$ws.Range("T14").Value = "Jämtland"
$ws.Range("U14").Value = "Krokom"
$ws.Range("V14").Value = "Jämtland"
$ws.Range("W14").Value = "Hotagen"
$ws.Range("Y14").Value = "'2023-09-05"
$ws.Range("Z14").Value = "00:00"
$ws.Range("AA14").Value = "'2023-09-05"
$ws.Range("AB14").Value = "00:00"
$ws.Range("AD14").Value = $false
$ws.Range("AE14").Value = $false
$ws.Range("AG14").Value = $false
$ws.Range("AT14").Value = "'"
$ws.Range("AW14").Value = "Benny Öwre"
$ws.Range("AX14").Value = "Benny Öwre"
$ws.Range("AY14").Value = "'"

# Row 15
$ws.Range("A15").Value = 111936795
$ws.Range("B15").Value = 56398
$ws.Range("C15").Value = "Ovaliderad"
$ws.Range("D15").Value = "NT"
$ws.Range("E15").Value = 100109
$ws.Range("F15").Value = "Tretåig hackspett"
$ws.Range("G15").Value = "Picoides tridactylus"
$ws.Range("H15").Value = "(Linnaeus, 1758)"
$ws.Range("I15").Value = "'"
$ws.Range("K15").Value = "'"
$ws.Range("L15").Value = "'"
$ws.Range("M15").Value = "'"
$ws.Range("N15").Value = "'"
$ws.Range("P15").Value = "rörvattsbodarna, Jmt"
$ws.Range("Q15").Value = 448749.3706757246
$ws.Range("R15").Value = 7087421.839990681
$ws.Range("S15").Value = 10
$ws.Range("T15").Value = "Jämtland"
$ws.Range("U15").Value = "Krokom"
$ws.Range("V15").Value = "Jämtland"
$ws.Range("W15").Value = "Hotagen"
$ws.Range("Y15").Value = "'2023-09-05"
$ws.Range("Z15").Value = "00:00"
$ws.Range("AA15").Value = "'2023-09-05"
$ws.Range("AB15").Value = "00:00"
$ws.Range("AC15").Value = "ringhack äldre"
$ws.Range("AD15").Value = $false
$ws.Range("AE15").Value = $false
$ws.Range("AG15").Value = $false
$ws.Range("AT15").Value = "'"
$ws.Range("AW15").Value = "Benny Öwre"
$ws.Range("AX15").Value = "Benny Öwre"
$ws.Range("AY15").Value = "'"

# Row 16
$ws.Range("A16").Value = 111936865
$ws.Range("B16").Value = 89423
$ws.Range("C16").Value = "Ovaliderad"
$ws.Range("D16").Value = "NT"
$ws.Range("E16").Value = 5432
$ws.Range("F16").Value = "Granticka"
$ws.Range("G16").Value = "Porodaedalea chrysoloma"
$ws.Range("H16").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("I16").Value = "'"
$ws.Range("P16").Value = "rörvattsbodarna, Jmt"
$ws.Range("Q16").Value = 448738.4239939091
$ws.Range("R16").Value = 7087426.42220111
$ws.Range("S16").Value = 10
$ws.Range("T16").Value = "Jämtland"
$ws.Range("U16").Value = "Krokom"
$ws.Range("V16").Value = "Jämtland"
$ws.Range("W16").Value = "Hotagen"
$ws.Range("Y16").Value = "'2023-09-05"
$ws.Range("Z16").Value = "00:00"
$ws.Range("AA16").Value = "'2023-09-05"
$ws.Range("AB16").Value = "00:00"
$ws.Range("AD16").Value = $false
$ws.Range("AE16").Value = $false
$ws.Range("AG16").Value = $false
$ws.Range("AT16").Value = "'"
$ws.Range("AW16").Value = "Benny Öwre"
$ws.Range("AX16").Value = "Benny Öwre"
$ws.Range("AY16").Value = "'"

# Row 17
$ws.Range("A17").Value = 111936796
$ws.Range("B17").Value = 56398
$ws.Range("C17").Value = "Ovaliderad"
$ws.Range("D17").Value = "NT"
$ws.Range("E17").Value = 100109
$ws.Range("F17").Value = "Tretåig hackspett"
$ws.Range("G17").Value = "Picoides tridactylus"
$ws.Range("H17").Value = "(Linnaeus, 1758)"
$ws.Range("I17").Value = "'"
$ws.Range("K17").Value = "'"
$ws.Range("L17").Value = "'"
$ws.Range("M17").Value = "'"
$ws.Range("N17").Value = "'"
$ws.Range("P17").Value = "rörvattsbodarna, Jmt"
$ws.Range("Q17").Value = 448882.8980770012
$ws.Range("R17").Value = 7087229.443335658
$ws.Range("S17").Value = 10
$ws.Range("T17").Value = "Jämtland"
$ws.Range("U17").Value = "Krokom"
$ws.Range("V17").Value = "Jämtland"
$ws.Range("W17").Value = "Hotagen"
$ws.Range("Y17").Value = "'2023-09-05"
$ws.Range("Z17").Value = "00:00"
$ws.Range("AA17").Value = "'2023-09-05"
$ws.Range("AB17").Value = "00:00"
$ws.Range("AC17").Value = "ringhack äldre"
$ws.Range("AD17").Value = $false
$ws.Range("AE17").Value = $false
$ws.Range("AG17").Value = $false
$ws.Range("AT17").Value = "'"
$ws.Range("AW17").Value = "Benny Öwre"
$ws.Range("AX17").Value = "Benny Öwre"
$ws.Range("AY17").Value = "'"

# Row 18
$ws.Range("A18").Value = 111936858
$ws.Range("B18").Value = 89845
$ws.Range("C18").Value = "Ovaliderad"
$ws.Range("D18").Value = "VU"
$ws.Range("E18").Value = 1209
$ws.Range("F18").Value = "Rynkskinn"
$ws.Range("G18").Value = "Phlebia centrifuga"
$ws.Range("H18").Value = "P.Karst."
$ws.Range("I18").Value = "'"
$ws.Range("P18").Value = "rörvattsbodarna, Jmt"
$ws.Range("Q18").Value = 448737.3665225056
$ws.Range("R18").Value = 7087496.445579056
$ws.Range("S18").Value = 10
$ws.Range("T18").Value = "Jämtland"
$ws.Range("U18").Value = "Krokom"
$ws.Range("V18").Value = "Jämtland"
$ws.Range("W18").Value = "Hotagen"
$ws.Range("Y18").Value = "'2023-09-05"
$ws.Range("Z18").Value = "00:00"
$ws.Range("AA18").Value = "'2023-09-05"
$ws.Range("AB18").Value = "00:00"
$ws.Range("AD18").Value = $false
$ws.Range("AE18").Value = $false
$ws.Range("AG18").Value = $false
$ws.Range("AT18").Value = "'"
$ws.Range("AW18").Value = "Benny Öwre"
$ws.Range("AX18").Value = "Benny Öwre"
$ws.Range("AY18").Value = "'"

# Row 19
$ws.Range("A19").Value = 111936893
$ws.Range("B19").Value = 77515
$ws.Range("C19").Value = "Ovaliderad"
$ws.Range("D19").Value = "NT"
$ws.Range("E19").Value = 6425
$ws.Range("F19").Value = "Garnlav"
$ws.Range("G19").Value = "Alectoria sarmentosa"
$ws.Range("H19").Value = "(Ach.) Ach."
$ws.Range("I19").Value = "'"
$ws.Range("P19").Value = "rörvattsbodarna, Jmt"
$ws.Range("Q19").Value = 448742.3011697636
$ws.Range("R19").Value = 7087501.648173723
$ws.Range("S19").Value = 10
$ws.Range("T19").Value = "Jämtland"
$ws.Range("U19").Value = "Krokom"
$ws.Range("V19").Value = "Jämtland"
$ws.Range("W19").Value = "Hotagen"
$ws.Range("Y19").Value = "'2023-09-05"
$ws.Range("Z19").Value = "00:00"
$ws.Range("AA19").Value = "'2023-09-05"
$ws.Range("AB19").Value = "00:00"
$ws.Range("AD19").Value = $false
$ws.Range("AE19").Value = $false
$ws.Range("AG19").Value = $false
$ws.Range("AT19").Value = "'"
$ws.Range("AW19").Value = "Benny Öwre"
$ws.Range("AX19").Value = "Benny Öwre"
$ws.Range("AY19").Value = "'"

# Row 20
$ws.Range("A20").Value = 111936870
$ws.Range("B20").Value = 89423
$ws.Range("C20").Value = "Ovaliderad"
$ws.Range("D20").Value = "NT"
$ws.Range("E20").Value = 5432
$ws.Range("F20").Value = "Granticka"
$ws.Range("G20").Value = "Porodaedalea chrysoloma"
$ws.Range("H20").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("I20").Value = "'"
$ws.Range("P20").Value = "rörvattsbodarna, Jmt"
$ws.Range("Q20").Value = 449019.027096529
$ws.Range("R20").Value = 7087276.979166135
$ws.Range("S20").Value = 10
$ws.Range("T20").Value = "Jämtland"
$ws.Range("U20").Value = "Krokom"
$ws.Range("V20").Value = "Jämtland"
$ws.Range("W20").Value = "Hotagen"
$ws.Range("Y20").Value = "'2023-09-05"
$ws.Range("Z20").Value = "00:00"
$ws.Range("AA20").Value = "'2023-09-05"
$ws.Range("AB20").Value = "00:00"
$ws.Range("AD20").Value = $false
$ws.Range("AE20").Value = $false
$ws.Range("AG20").Value = $false
$ws.Range("AT20").Value = "'"
$ws.Range("AW20").Value = "Benny Öwre"
$ws.Range("AX20").Value = "Benny Öwre"
$ws.Range("AY20").Value = "'"

# Row 21
$ws.Range("A21").Value = 111936869
$ws.Range("B21").Value = 89423
$ws.Range("C21").Value = "Ovaliderad"
$ws.Range("D21").Value = "NT"
$ws.Range("E21").Value = 5432
$ws.Range("F21").Value = "Granticka"
$ws.Range("G21").Value = "Porodaedalea chrysoloma"
$ws.Range("H21").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("I21").Value = "'"
$ws.Range("P21").Value = "rörvattsbodarna, Jmt"
$ws.Range("Q21").Value = 449143.8568242944
$ws.Range("R21").Value = 7087117.752608996
$ws.Range("S21").Value = 10
$ws.Range("T21").Value = "Jämtland"
$ws.Range("U21").Value = "Krokom"
$ws.Range("V21").Value = "Jämtland"
$ws.Range("W21").Value = "Hotagen"
$ws.Range("Y21").Value = "'2023-09-05"
$ws.Range("Z21").Value = "00:00"
$ws.Range("AA21").Value = "'2023-09-05"
$ws.Range("AB21").Value = "00:00"
$ws.Range("AD21").Value = $false
$ws.Range("AE21").Value = $false
$ws.Range("AG21").Value = $false
$ws.Range("AT21").Value = "'"
$ws.Range("AW21").Value = "Benny Öwre"
$ws.Range("AX21").Value = "Benny Öwre"
$ws.Range("AY21").Value = "'"

